$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.608580112457275
$ws.Range("B1").Value = 2.45148491859436
$ws.Range("C1").Value = 1.843499302864075
$ws.Range("D1").Value = 1.74121356010437
$ws.Range("E1").Value = 1.836049199104309
